# HOTFIX change date format in example
# Change the date-like text values in columns K:P (Дата поверки, Дата следующей
# поверки, Дата установки, Дата ввода в эксплуатацию, Дата опломбирования,
# Дата контрольных показаний) from DD.MM.YYYY to ISO YYYY-MM-DD format,
# for every data row (2 through 10), unifying them to a single sequential
# run of dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("2021-01-20", "2021-01-21", "2021-01-22", "2021-01-23", "2021-01-24", "2021-01-25")
$columns = @("K", "L", "M", "N", "O", "P")

for ($row = 2; $row -le 10; $row++) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $dates[$i]
    }
}
